# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Columna A (grupo-cultivo-descripcion): pasa de dimension a measure
# Columna E (municipio-nombre): pasa de measure a dimension (sdmx refArea)
# Columna H (cultivo-detalle-descripcion): pasa de dimension a measure
# Las filas de mapping (A5/H5) de las columnas que ahora son measure se eliminan

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columna A: grupo-cultivo-descripcion (dimension -> measure) ---
$ws.Range("A2").Value = "iaest-measure:grupo-cultivo-descripcion"
$ws.Range("A3").Value = "medida"
$ws.Range("A4").Value = "xsd:int"

# --- Columna E: municipio-nombre (measure -> dimension, sdmx refArea) ---
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"

# --- Columna H: cultivo-detalle-descripcion (dimension -> measure) ---
$ws.Range("H2").Value = "iaest-measure:cultivo-detalle-descripcion"
$ws.Range("H3").Value = "medida"
$ws.Range("H4").Value = "xsd:int"

# --- Fila 5: ya no hay fichero de mapping para columnas A y H (ahora measure) ---
$ws.Range("A5").Clear()
$ws.Range("H5").Clear()
